# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> stock "Office Theme" palette (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" palette (used by the Slide Master / all slides)
#
# The authored change swaps the two themes' content so the slides (which are
# driven by the Slide Master's theme) now render with the default "Office
# Theme" colors instead of "Integral". The font scheme and format scheme
# (fills/lines/effects) are already identical between the two themes, so the
# only substantive difference is the color scheme (clrScheme) - this script
# re-points every theme color slot used by the presentation's (one and only)
# design to the stock Office Theme RGB values.
#
# PowerPoint's ColorFormat.RGB is packed 0x00BBGGRR (blue/green/red), so each
# literal below is the target hex color (RRGGBB) with its bytes reversed.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$cs = $s1.ThemeColorScheme

# clrScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$cs.Colors(1).RGB  = 0x000000   # dk1      -> RGB 000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> RGB FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      -> RGB 44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> RGB E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  -> RGB 5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  -> RGB ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> RGB A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  -> RGB FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  -> RGB 4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  -> RGB 70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    -> RGB 0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink -> RGB 954F72
